# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Freshly-scraped K (strikeout) values per game, replacing the old
# "Strike#" figures that used to live in column G ("K" header).
$kVals = @(
    1, 3, 3, 2, 1, 1, 3, 1, 0, 2,
    2, 1, 0, 2, 2, 1, 2, 2, 3, 2,
    1, 1, 2, 0, 0, 1, 2, 1, 3, 0,
    2, 1, 0, 3, 2, 2, 3, 2, 0, 2,
    3, 1, 2, 0, 0, 1, 3, 2, 0, 2,
    2, 2, 2, 1, 2, 0, 2, 1, 3, 2,
    2, 1
)

$firstRow = 2
for ($i = 0; $i -lt $kVals.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kVals[$i]
}
